# Updated cryptos list on Fri Apr 21 09:36:04 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings stored as text in the
# original file (e.g. "1.003", "28.095.16"). Force each cell to the
# Text number format before assignment so Excel does not auto-convert
# these strings into Number values.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.095.16"
$ws.Range("E2").Value = "  -3.44%  "
$ws.Range("D3").Value = "1.925.33"
$ws.Range("E3").Value = "  -2.36%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "328.95"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").Value = "0.4725"
$ws.Range("E7").Value = "  -5.30%  "
$ws.Range("D8").Value = "0.4067"
$ws.Range("E8").Value = "  -3.60%  "
$ws.Range("D9").Value = "53.04"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "0.08440"
$ws.Range("E10").Value = "  -8.82%  "
$ws.Range("E11").Value = "  -4.90%  "
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").Value = "1.932.34"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").Value = "7.513"
$ws.Range("D15").Value = "6.101"
$ws.Range("E15").Value = "  -5.53%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  -0.55%  "
$ws.Range("D17").Value = "90.83"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "0.00001067"
$ws.Range("E18").Value = "  -3.78%  "
$ws.Range("D19").Value = "0.06587"
$ws.Range("E19").Value = "  -2.09%  "
$ws.Range("D20").Value = "18.08"
$ws.Range("E20").Value = "  -6.17%  "
$ws.Range("D21").Value = "1.003"
$ws.Range("D22").Value = "5.760"
$ws.Range("E22").Value = "  -3.54%  "
$ws.Range("D23").Value = "28.101.67"
$ws.Range("E23").Value = "  -3.50%  "
$ws.Range("E24").Value = "  -4.48%  "
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "2.167.49"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("D27").Value = "154.37"
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").Value = "20.13"
$ws.Range("E28").Value = "  -2.98%  "
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("D30").Value = "5.721"
$ws.Range("E30").Value = "  -9.33%  "
$ws.Range("D31").Value = "123.81"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "0.9748"
$ws.Range("E32").Value = "  -7.28%  "
$ws.Range("D33").Value = "0.09600"
$ws.Range("E33").Value = "  -2.70%  "
$ws.Range("D34").Value = "1.442"
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "3.641"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "5.559"
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("D37").Value = "9.041"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "0.02318"
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("D39").Value = "0.06166"
$ws.Range("E39").Value = "  -3.89%  "
$ws.Range("D40").Value = "1.237"
$ws.Range("E40").Value = "  -6.35%  "
$ws.Range("D41").Value = "0.6171"
$ws.Range("E41").Value = "  -4.83%  "
$ws.Range("D42").Value = "11.07"
$ws.Range("E42").Value = "  -3.73%  "
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "0.1907"
$ws.Range("E44").Value = "  -4.88%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5899"
$ws.Range("E45").Value = "  -5.27%  "
$ws.Range("B46").Value = "WEMIXTOKEN"
$ws.Range("C46").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D46").Value = "1.291"
$ws.Range("E46").Value = "  -6.86%  "
$ws.Range("D47").Value = "12.84"
$ws.Range("E47").Value = "  -4.11%  "
$ws.Range("D48").Value = "2.037"
$ws.Range("E48").Value = "  -7.13%  "
$ws.Range("D49").Value = "3.468"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "0.06819"
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").Value = "110.05"
$ws.Range("E51").Value = "  -2.96%  "
